$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1231
$ws1.Range("F4").Value = 12798
$ws1.Range("F5").Value = 726
$ws1.Range("F7").Value = 318
$ws1.Range("F10").Value = 1866
$ws1.Range("F14").Value = 208
$ws1.Range("F16").Value = 344
$ws1.Range("F17").Value = 228
$ws1.Range("F18").Value = 295
$ws1.Range("F20").Value = 126
$ws1.Range("F22").Value = 217
$ws1.Range("F23").Value = 245
$ws1.Range("F24").Value = 1285
$ws1.Range("F25").Value = 334

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 282
$ws2.Range("F5").Value = 4452
$ws2.Range("F6").Value = 155
$ws2.Range("F8").Value = 11
$ws2.Range("F11").Value = 357
$ws2.Range("F17").Value = 11

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 872
$ws3.Range("F3").Value = 4004
$ws3.Range("G3").Value = "已售罄"

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 872
$ws4.Range("F6").Value = 1231
$ws4.Range("F7").Value = 12798
$ws4.Range("F8").Value = 282
$ws4.Range("F9").Value = 726
$ws4.Range("F10").Value = 4004
$ws4.Range("G10").Value = "已售罄"
$ws4.Range("F12").Value = 318
$ws4.Range("F15").Value = 1866
$ws4.Range("F19").Value = 4452
$ws4.Range("F20").Value = 208
$ws4.Range("F21").Value = 155
$ws4.Range("F22").Value = 155
$ws4.Range("F25").Value = 11
$ws4.Range("F28").Value = 357
$ws4.Range("F29").Value = 344
$ws4.Range("F31").Value = 228
$ws4.Range("F32").Value = 295
$ws4.Range("F34").Value = 126
$ws4.Range("F37").Value = 217
$ws4.Range("F40").Value = 245
$ws4.Range("F41").Value = 1285
$ws4.Range("F43").Value = 334
$ws4.Range("F46").Value = 11
